# Resume formatting fix:
#   1. Split the run "15-20 students), Technical work with sound equipment."
#      into separate runs (same formatting) at word/phrase boundaries.
#   2. Insert a paragraph break right after "sound equipment." so the
#      trailing tab (and the "_GoBack" bookmark) move into a new paragraph
#      of their own (reusing the same tab stops / run formatting).
#   3. Relocate the "_GoBack" bookmark to sit at the start of that new
#      paragraph (it used to sit just before "15-20").
#   4. Shrink the page bottom margin from 1440 twips (72pt) to 1080 twips
#      (54pt).

$d = $word.ActiveDocument

$target = "15-20 students), Technical work with sound equipment."
$rng = $d.Content
[void]$rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $rng.Start

# Relative offsets of the new run boundaries within $target:
#   "15-20 "                     0-6
#   "students), Technical work"  6-31
#   " "                          31-32
#   "with "                      32-37
#   "sound equipment."           37-53
# Toggling Bold on then back off over each sub-range forces the engine to
# split out a distinct run there without leaving any visible/formatting
# change behind.
$splits = @(0, 6, 31, 32, 37)
for ($i = 0; $i -lt $splits.Length - 1; $i++) {
    $sub = $d.Range($base + $splits[$i], $base + $splits[$i + 1])
    $sub.Bold = 1
    $sub.Bold = 0
}

# Break the paragraph right after "sound equipment." (offset 53), pushing
# the bookmark + trailing tab run into a new paragraph that inherits the
# same paragraph formatting (tab stops, run properties).
$breakPoint = $d.Range($base + 53, $base + 53)
[void]$breakPoint.InsertParagraphAfter()

# Move the "_GoBack" bookmark so it starts the new paragraph (adding a
# bookmark with a name that already exists relocates it).
$newParaStart = $d.Range($base + 54, $base + 54)
$d.Bookmarks.Add("_GoBack", $newParaStart)

# Bottom margin: 1440 twips -> 1080 twips (PageSetup measurements are in
# points; 20 twips = 1 point, so 1080 twips = 54pt).
$d.Sections.Item(1).PageSetup.BottomMargin = 54
